$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Range("F1").Value = "First day - ENEM 2019"
$ws.Range("G1").Value = "Second day - ENEM 2019"

$values = @{
    2  = @(14.43, 18.93)
    3  = @(9.32, 13.81)
    4  = @(16.49, 21.57)
    5  = @(13.99, 20.35)
    6  = @(10.26, 13.75)
    7  = @(10.58, 15.08)
    8  = @(11.14, 14.75)
    9  = @(12.91, 16.05)
    10 = @(9.81, 11.92)
    11 = @(17.9, 20.92)
    12 = @(8.869999999999999, 11.58)
    13 = @(10.08, 12.69)
    14 = @(13.2, 16.51)
    15 = @(11.69, 14.86)
    16 = @(10.15, 12.86)
    17 = @(11.69, 15.13)
    18 = @(11.42, 15)
    19 = @(10.11, 14.17)
    20 = @(11.11, 15.28)
    21 = @(14.2, 19.16)
    22 = @(14.25, 19.43)
    23 = @(13.47, 18.38)
    24 = @(10.82, 15.34)
    25 = @(10.66, 15.14)
    26 = @(11.94, 16.51)
    27 = @(10.41, 14.05)
    28 = @(8.17, 11.29)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 6).Value = $pair[0]
    $ws.Cells.Item($row, 7).Value = $pair[1]
}
